$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Periodo Mora" values in column E (rows 16-19) are reordered so that the
# most recent period comes first: 2305,2306,2307,2308 -> 2308,2307,2306,2305
$ws.Range("E16").Value = "2308"
$ws.Range("E17").Value = "2307"
$ws.Range("E18").Value = "2306"
$ws.Range("E19").Value = "2305"
